$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31 (shifts existing rows 31..107 down to 32..108),
# carrying the existing column formatting (e.g. the date style on column D) along.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with a new daily price record for
# Ajo / Chino / Primera at Vega Monumental Concepción, matching the
# other rows' fixed fields and the new date/volume/price figures.
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = 'Vega Monumental Concepción'
$ws.Range("C31").Value = 'Bíobío'
$ws.Range("D31").Value = 44519
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100112003
$ws.Range("G31").Value = 'Ajo'
$ws.Range("H31").Value = 'Chino'
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 250
$ws.Range("K31").Value = 16000
$ws.Range("L31").Value = 17000
$ws.Range("M31").Value = 16400
$ws.Range("N31").Value = '$/caja 10 kilos'
$ws.Range("O31").Value = 'China'
$ws.Range("P31").Value = 1640
$ws.Range("Q31").Value = 10
$ws.Range("R31").Value = 'Hortaliza'
